$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (the source data is text, e.g. "26.20").
# NumberFormat must be set per-cell: multi-area ranges ("D6,D7") only honour
# the first area for this property.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated price / 1h-volume figures scraped by the Action run.
$ws.Range("D2").Value = '67.921.47'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '3.470.01'
$ws.Range("E3").Value = '  -0.71%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '592.34'
$ws.Range("E5").Value = '  -1.08%  '

$ws.Range("D6").Value = '181.33'
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +3.65%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = '3.467.44'
$ws.Range("E9").Value = '  -0.74%  '

$ws.Range("D10").Value = '0.140'
$ws.Range("E10").Value = '  +0.42%  '

$ws.Range("D11").Value = '6.99'
$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("D12").Value = '0.430'
$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("D13").Value = '4.071.97'
$ws.Range("E13").Value = '  -0.83%  '

$ws.Range("D14").Value = '32.13'
$ws.Range("E14").Value = '  -0.66%  '

$ws.Range("D16").Value = '67.897.50'
$ws.Range("E16").Value = '  +0.36%  '

$ws.Range("D17").Value = '0.0000177'
$ws.Range("E17").Value = '  -2.15%  '

$ws.Range("D18").Value = '3.483.89'
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").Value = '6.21'
$ws.Range("E19").Value = '  -2.36%  '

$ws.Range("D20").Value = '14.07'
$ws.Range("E20").Value = '  -4.16%  '

$ws.Range("D21").Value = '393.11'
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").Value = '7.92'
$ws.Range("E22").Value = '  -1.66%  '

$ws.Range("E23").Value = '  +2.33%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = '0.538'
$ws.Range("E25").Value = '  -1.05%  '

$ws.Range("D26").Value = '71.94'
$ws.Range("E26").Value = '  -1.97%  '

$ws.Range("D27").Value = '0.0000120'
$ws.Range("E27").Value = '  -3.43%  '

$ws.Range("D28").Value = '10.41'
$ws.Range("E28").Value = '  -0.12%  '

$ws.Range("D29").Value = '0.176'
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("D31").Value = '6.12'
$ws.Range("E31").Value = '  -2.46%  '

$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -0.53%  '

$ws.Range("D33").Value = '1.40'
$ws.Range("E33").Value = '  -4.31%  '

$ws.Range("D34").Value = '23.59'
$ws.Range("E34").Value = '  -1.50%  '

$ws.Range("D35").Value = '7.36'
$ws.Range("E35").Value = '  -0.94%  '

$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").Value = '1.56'
$ws.Range("E37").Value = '  -6.11%  '

$ws.Range("D38").Value = '161.70'
$ws.Range("E38").Value = '  -1.58%  '

$ws.Range("D39").Value = '0.887'
$ws.Range("E39").Value = '  +1.60%  '

$ws.Range("D40").Value = '2.81'
$ws.Range("E40").Value = '  +4.38%  '

$ws.Range("D41").Value = '1.87'
$ws.Range("E41").Value = '  -4.26%  '

$ws.Range("D42").Value = '4.65'
$ws.Range("E42").Value = '  -1.60%  '

$ws.Range("D43").Value = '6.72'
$ws.Range("E43").Value = '  -5.40%  '

$ws.Range("D44").Value = '26.07'
$ws.Range("E44").Value = '  -2.00%  '

$ws.Range("D45").Value = '0.0718'
$ws.Range("E45").Value = '  -2.46%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '26.20'
$ws.Range("E46").Value = '  -5.41%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.732.31'
$ws.Range("E47").Value = '  -2.96%  '

$ws.Range("D48").Value = '41.34'
$ws.Range("E48").Value = '  -2.35%  '

$ws.Range("D49").Value = '0.0298'
$ws.Range("E49").Value = '  -2.13%  '

$ws.Range("D50").Value = '327.53'
$ws.Range("E50").Value = '  -5.10%  '

$ws.Range("D51").Value = '1.04'
$ws.Range("E51").Value = '  -3.56%  '
